$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Fill in evaluation scores (value 1) for all evaluators in rows 5-9, columns B-E
# (columns F and G already contain 1 in the source file)
$ws.Range("B5:E9").Value = 1

# Mark the free-text "other characteristics" row (row 10) as not applicable
$ws.Range("B10:G10").Value = "-"

# Keep the last active selection consistent with the authored workbook
$ws.Range("H11").Select() | Out-Null
